# Apply the "Optuna Attempt (go back with original)" edits to the workbook.
# Updates numeric metrics on the "Forecast Comparison" sheet and the
# summary totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: numeric cell updates -----------------
$forecastUpdates = @{
    "L2"  = 1.18
    "D3"  = 43
    "H3"  = 2.54
    "L3"  = 1.09
    "H4"  = 1.89
    "L4"  = 0.9
    "H5"  = 0.8
    "L5"  = 1.17
    "L6"  = 0.82
    "L7"  = 0.84
    "L8"  = 1.19
    "L9"  = 0.92
    "L10" = 1.1
    "D11" = 40
    "L11" = 0.8
    "L12" = 0.86
    "D13" = 39
    "L13" = 0.99
    "D14" = 39
    "L14" = 0.99
    "L15" = 0.89
    "L16" = 1.06
    "L17" = 1.1
}

foreach ($addr in $forecastUpdates.Keys) {
    $wsForecast.Range($addr).Value = $forecastUpdates[$addr]
}

# --- Summary sheet: totals stored as text values ----------------------
# These cells are text (not numeric) in the workbook, so a leading
# apostrophe keeps Excel from auto-converting them to numbers, and
# resetting the Style afterward avoids leaving a "quote prefix" style
# applied to the cell.
$summaryUpdates = @{
    "B9"  = "619"
    "B10" = "313"
    "B11" = "159"
    "B12" = "43"
}

foreach ($addr in $summaryUpdates.Keys) {
    $cell = $wsSummary.Range($addr)
    $cell.Value = "'" + $summaryUpdates[$addr]
    $cell.Style = "Normal"
}
